$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'57905035"
$ws.Range("C2").Value = "FCT903349878300606464"
$ws.Range("F2").Value = "'697.36"
$ws.Range("I2").Value = "57905035+1"

$ws.Range("E2").Select()
